{"js": "// Replace the date line and the 25 multiplication problems with their new\n// values. Each original string is unique in the document, so a simple\n// search-and-replace per pair is unambiguous and order-independent.\nconst replacements = [\n  [\"2025-06-21 Saturday\", \"2025-06-22 Sunday\"],\n  [\"337\u00d74=\", \"818\u00d75=\"],\n  [\"991\u00d79=\", \"773\u00d75=\"],\n  [\"430\u00d79=\", \"354\u00d76=\"],\n  [\"378\u00d75=\", \"304\u00d79=\"],\n  [\"946\u00d78=\", \"531\u00d73=\"],\n  [\"893\u00d78=\", \"855\u00d73=\"],\n  [\"691\u00d78=\", \"577\u00d78=\"],\n  [\"409\u00d77=\", \"585\u00d72=\"],\n  [\"726\u00d74=\", \"976\u00d74=\"],\n  [\"825\u00d74=\", \"313\u00d78=\"],\n  [\"301\u00d76=\", \"908\u00d72=\"],\n  [\"170\u00d78=\", \"995\u00d78=\"],\n  [\"894\u00d79=\", \"594\u00d75=\"],\n  [\"980\u00d74=\", \"975\u00d72=\"],\n  [\"525\u00d79=\", \"796\u00d73=\"],\n  [\"435\u00d78=\", \"280\u00d77=\"],\n  [\"284\u00d74=\", \"134\u00d73=\"],\n  [\"173\u00d73=\", \"794\u00d73=\"],\n  [\"147\u00d74=\", \"538\u00d77=\"],\n  [\"449\u00d75=\", \"355\u00d79=\"],\n  [\"387\u00d77=\", \"728\u00d76=\"],\n  [\"923\u00d75=\", \"980\u00d78=\"],\n  [\"274\u00d76=\", \"588\u00d79=\"],\n  [\"516\u00d76=\", \"964\u00d77=\"],\n  [\"244\u00d76=\", \"219\u00d73=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the 25 multiplication problems with their new\n# values. Each original string is unique in the document, so a simple\n# Find/Replace (ReplaceAll) per pair is unambiguous and order-independent.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-06-21 Saturday\", \"2025-06-22 Sunday\"),\n    @(\"337\u00d74=\", \"818\u00d75=\"),\n    @(\"991\u00d79=\", \"773\u00d75=\"),\n    @(\"430\u00d79=\", \"354\u00d76=\"),\n    @(\"378\u00d75=\", \"304\u00d79=\"),\n    @(\"946\u00d78=\", \"531\u00d73=\"),\n    @(\"893\u00d78=\", \"855\u00d73=\"),\n    @(\"691\u00d78=\", \"577\u00d78=\"),\n    @(\"409\u00d77=\", \"585\u00d72=\"),\n    @(\"726\u00d74=\", \"976\u00d74=\"),\n    @(\"825\u00d74=\", \"313\u00d78=\"),\n    @(\"301\u00d76=\", \"908\u00d72=\"),\n    @(\"170\u00d78=\", \"995\u00d78=\"),\n    @(\"894\u00d79=\", \"594\u00d75=\"),\n    @(\"980\u00d74=\", \"975\u00d72=\"),\n    @(\"525\u00d79=\", \"796\u00d73=\"),\n    @(\"435\u00d78=\", \"280\u00d77=\"),\n    @(\"284\u00d74=\", \"134\u00d73=\"),\n    @(\"173\u00d73=\", \"794\u00d73=\"),\n    @(\"147\u00d74=\", \"538\u00d77=\"),\n    @(\"449\u00d75=\", \"355\u00d79=\"),\n    @(\"387\u00d77=\", \"728\u00d76=\"),\n    @(\"923\u00d75=\", \"980\u00d78=\"),\n    @(\"274\u00d76=\", \"588\u00d79=\"),\n    @(\"516\u00d76=\", \"964\u00d77=\"),\n    @(\"244\u00d76=\", \"219\u00d73=\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $pair[0]\n    $replace = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Text = $find\n    $rng.Find.Replacement.Text = $replace\n    $rng.Find.Forward = $true\n    $rng.Find.Wrap = 1  # wdFindContinue\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWholeWord = $false\n    $rng.Find.MatchWildcards = $false\n    $rng.Find.Execute([ref]$find, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$replace, [ref]2) | Out-Null\n}\n"}
